# Canada Premier League workbook update (14-04-2024 15:12)
#
# 1) Rows 11 and 12 (match ids 9 and 10) had their data swapped back —
#    every column except A (the running "id" column) is exchanged between
#    the two rows.
# 2) A brand-new fixture (row 89, id 87) is appended at the bottom of the
#    sheet for an upcoming (not-yet-played) Forge FC vs Cavalry FC match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Step 1: swap the contents (values + formats) of B11:AC11 with B12:AC12,
# using a scratch row far below the data as a temporary holding area,
# then deleting that scratch row entirely so it leaves no trace behind.
# ---------------------------------------------------------------------
$row11 = $ws.Range("B11:AC11")
$row12 = $ws.Range("B12:AC12")
$tmp = $ws.Range("B200:AC200")

$row11.Copy($tmp)
$row12.Copy($row11)
$tmp.Copy($row12)
$ws.Rows.Item(200).Delete()

# ---------------------------------------------------------------------
# Step 2: append the new fixture as row 89 (copy row 88's formatting,
# then overwrite with the new fixture's values).
# ---------------------------------------------------------------------
$srcRow = $ws.Range("A88:AC88")
$dstRow = $ws.Range("A89:AC89")
$srcRow.Copy($dstRow)

# This fixture has no result yet, so clear the result-derived columns
# (FTHG, FTAG, FTR, PL_AhOver, PL_AhUnder) that the copy brought along.
$ws.Range("H89:J89").ClearContents()
$ws.Range("AB89:AC89").ClearContents()

$ws.Cells.Item(89, 1).Value = 87
$ws.Cells.Item(89, 2).Value = 7802874
$ws.Cells.Item(89, 3).Value = "Canada Premier League"
$ws.Cells.Item(89, 4).Value = "Canada Premier League"
$ws.Cells.Item(89, 5).Value = 45395.70833333334
$ws.Cells.Item(89, 6).Value = "Forge FC"
$ws.Cells.Item(89, 7).Value = "Cavalry FC"

$ws.Cells.Item(89, 11).Value = 2.2
$ws.Cells.Item(89, 12).Value = 3.6
$ws.Cells.Item(89, 13).Value = 2.6
$ws.Cells.Item(89, 14).Value = 2.25
$ws.Cells.Item(89, 15).Value = 3.6
$ws.Cells.Item(89, 16).Value = 2.55
$ws.Cells.Item(89, 17).Value = 0
$ws.Cells.Item(89, 18).Value = 1.775
$ws.Cells.Item(89, 19).Value = 2.025
$ws.Cells.Item(89, 20).Value = 2.5
$ws.Cells.Item(89, 21).Value = 1.95
$ws.Cells.Item(89, 22).Value = 1.85
$ws.Cells.Item(89, 23).Value = 0
$ws.Cells.Item(89, 24).Value = 0
$ws.Cells.Item(89, 25).Value = 0
$ws.Cells.Item(89, 26).Value = 0
$ws.Cells.Item(89, 27).Value = 0
